# Auto-generated PowerShell COM-interop script applying the diff changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 394.6875
$ws.Range("I2").Value = 296.5
$ws.Range("K2").Value = 296.5
$ws.Range("M2").Value = -183.5
$ws.Range("H28").Value = 285.72726
$ws.Range("I28").Value = 229.22223
$ws.Range("K28").Value = 229.22223
$ws.Range("M28").Value = 255.77777
$ws.Range("H40").Value = 3120
$ws.Range("J40").Value = 3240
$ws.Range("L40").Value = 3240
$ws.Range("N40").Value = -3590
$ws.Range("H107").Value = 2498.25
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
$ws.Range("H109").Value = 68000
$ws.Range("J109").Value = 68000
$ws.Range("L109").Value = 68000
$ws.Range("N109").Value = -70774
$ws.Range("H111").Value = 1224.75
$ws.Range("I111").Value = 633
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 1899
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = 1168
$ws.Range("N111").Value = -15134
$ws.Range("H132").Value = 833.8788
$ws.Range("I132").Value = 775.918
$ws.Range("K132").Value = 2327.754
$ws.Range("M132").Value = 202.2460000000001
$ws.Range("H137").Value = 53981.684
$ws.Range("I137").Value = 954.4
$ws.Range("J137").Value = 112900.89
$ws.Range("K137").Value = 2863.2
$ws.Range("L137").Value = 338702.67
$ws.Range("M137").Value = -313.1999999999998
$ws.Range("N137").Value = -343802.67
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 242393.27
$ws.Range("I2").Value = 309235.88
$ws.Range("J2").Value = 1759.8
$ws.Range("K2").Value = 309235.88
$ws.Range("L2").Value = 1759.8
$ws.Range("M2").Value = -309122.88
$ws.Range("N2").Value = -1985.8
$ws.Range("H26").Value = 25000
$ws.Range("I26").Value = 25000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 25000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -24670
$ws.Range("N26").ClearContents()
$ws.Range("H32").Value = 5596.575
$ws.Range("I32").Value = 3438.1614
$ws.Range("K32").Value = 3438.1614
$ws.Range("M32").Value = -3151.1614
$ws.Range("H74").Value = 776.5
$ws.Range("I74").Value = 562.94116
$ws.Range("J74").Value = 4407
$ws.Range("K74").Value = 562.94116
$ws.Range("L74").Value = 4407
$ws.Range("M74").Value = 311.05884
$ws.Range("N74").Value = -6155
$ws.Range("H77").Value = 776.5
$ws.Range("I77").Value = 562.94116
$ws.Range("J77").Value = 4407
$ws.Range("K77").Value = 2814.7058
$ws.Range("L77").Value = 22035
$ws.Range("M77").Value = 1553.2942
$ws.Range("N77").Value = -30771
$ws.Range("H116").Value = 242393.27
$ws.Range("I116").Value = 309235.88
$ws.Range("J116").Value = 1759.8
$ws.Range("K116").Value = 309235.88
$ws.Range("L116").Value = 1759.8
$ws.Range("M116").Value = -306941.88
$ws.Range("N116").Value = -6347.8
$ws.Range("H132").Value = 1612.986
$ws.Range("I132").Value = 1304.9269
$ws.Range("J132").Value = 2034
$ws.Range("K132").Value = 3914.7807
$ws.Range("L132").Value = 6102
$ws.Range("M132").Value = -1384.7807
$ws.Range("N132").Value = -11162
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 242393.27
$ws.Range("I3").Value = 309235.88
$ws.Range("J3").Value = 1759.8
$ws.Range("K3").Value = 309235.88
$ws.Range("L3").Value = 1759.8
$ws.Range("M3").Value = -309121.88
$ws.Range("N3").Value = -1987.8
$ws.Range("H134").Value = 4027.6492
$ws.Range("I134").Value = 4411.7075
$ws.Range("K134").Value = 13235.1225
$ws.Range("M134").Value = -10700.1225
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 100010
$ws.Range("J23").Value = 100010
$ws.Range("L23").Value = 100010
$ws.Range("N23").Value = -100490
$ws.Range("H27").Value = 100010
$ws.Range("J27").Value = 100010
$ws.Range("L27").Value = 100010
$ws.Range("N27").Value = -100394
$ws.Range("H31").Value = 2380.8572
$ws.Range("J31").Value = 3001.4285
$ws.Range("L31").Value = 3001.4285
$ws.Range("N31").Value = -3591.4285
$ws.Range("H34").Value = 2380.8572
$ws.Range("J34").Value = 3001.4285
$ws.Range("L34").Value = 3001.4285
$ws.Range("N34").Value = -3405.4285
$ws.Range("H134").Value = 1809.7142
$ws.Range("I134").Value = 1528.5883
$ws.Range("K134").Value = 4585.7649
$ws.Range("M134").Value = -2050.7649
$ws.Range("H141").Value = 55083.93
$ws.Range("J141").Value = 53782.69
$ws.Range("L141").Value = 53782.69
$ws.Range("N141").Value = -64142.69
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 227.4
$ws.Range("J2").Value = 188.16667
$ws.Range("L2").Value = 1129.00002
$ws.Range("N2").Value = -1355.00002
$ws.Range("H33").Value = 230.54546
$ws.Range("I33").Value = 121.28571
$ws.Range("J33").Value = 421.75
$ws.Range("K33").Value = 727.71426
$ws.Range("L33").Value = 2530.5
$ws.Range("M33").Value = -444.71426
$ws.Range("N33").Value = -3096.5
$ws.Range("H38").Value = 380.81818
$ws.Range("J38").Value = 800.5
$ws.Range("L38").Value = 2401.5
$ws.Range("N38").Value = -3095.5
$ws.Range("H51").Value = 1500
$ws.Range("I51").Value = 1500
$ws.Range("K51").Value = 4500
$ws.Range("M51").Value = -4040
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 20285
$ws.Range("J109").Value = 20285
$ws.Range("L109").Value = 20285
$ws.Range("N109").Value = -22365
$ws.Range("H113").Value = 928.5
$ws.Range("I113").Value = 699
$ws.Range("K113").Value = 699
$ws.Range("M113").Value = 1471
$ws.Range("H132").Value = 917843.7
$ws.Range("I132").Value = 1167140.1
$ws.Range("K132").Value = 3501420.3
$ws.Range("M132").Value = -3498890.3
$ws.Range("H136").Value = 5663.231
$ws.Range("J136").Value = 5663.231
$ws.Range("L136").Value = 16989.693
$ws.Range("N136").Value = -22089.693
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 525.9524
$ws.Range("J55").Value = 678.5
$ws.Range("L55").Value = 678.5
$ws.Range("N55").Value = -1024.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 14799.8
$ws.Range("J21").Value = 14799.8
$ws.Range("L21").Value = 14799.8
$ws.Range("N21").Value = -15269.8
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H35").Value = 14799.8
$ws.Range("J35").Value = 14799.8
$ws.Range("L35").Value = 14799.8
$ws.Range("N35").Value = -15379.8
$ws.Range("H96").Value = 5401.95
$ws.Range("I96").Value = 1756.2222
$ws.Range("K96").Value = 1756.2222
$ws.Range("M96").Value = -383.2221999999999
$ws.Range("H122").Value = 35310.5
$ws.Range("J122").Value = 2255
$ws.Range("L122").Value = 6765
$ws.Range("N122").Value = -11665
